# Added verification point on flight search result page.
# Updates the "from"/"to" result cells on Sheet2 with the full
# flight-search-result strings shown on the results page.

$wb = $excel.ActiveWorkbook

$ws2 = $wb.Worksheets.Item("Sheet2")

$ws2.Range("A2").Value = "Chicago, IL, US (CHI - All Airports)"
$ws2.Range("B2").Value = "New York, NY, US (NYC - All Airports)"
$ws2.Range("A3").Value = "Washington, DC, US (IAD - Dulles)"
$ws2.Range("B3").Value = "San Francisco, CA, US (SFO)"

$ws2.Activate()
$ws2.Range("B10").Select()
